# "High - Low Heat" edit
# Splits the generic "heat" node into "heat_high" / "heat_low", wired
# together through a new "Auxilliary" unit ("heat_split") that sits on
# the Units sheet, and updates the Connections sheet's heat pipeline to
# reference the new "heat_low" node.

$wb = $excel.ActiveWorkbook
$units = $wb.Worksheets.Item("Units")
$connections = $wb.Worksheets.Item("Connections")

# --- Units sheet --------------------------------------------------------
# New unit: the "Auxilliary" heat_split unit, added below the table as row 9.
$units.Range("A9").Value2 = "heat_split"
$units.Range("B9").Value2 = "Auxilliary"
$units.Range("C9").Value2 = "heat_high"
$units.Range("E9").Value2 = "internal_heat"
$units.Range("F9").Value2 = "heat_low"
$units.Range("U9").Value2 = 0.4

# Electrolyzer (row 3) byproduct output becomes low-grade heat.
$units.Range("F3").Value2 = "heat_low"

# Fuel synthesizer (row 7) byproduct output becomes high-grade heat.
$units.Range("F7").Value2 = "heat_high"

# Match the formatting used by the other "resolution_output"/"demand"
# dropdown cells (right aligned) for the two new cells beyond the table.
$units.Range("AI9").HorizontalAlignment = -4152
$units.Range("AJ9").HorizontalAlignment = -4152

# Extend the "h, D, W, M, Q, Y" list validation to the new AI9 cell.
$units.Range("AI9").Validation.Add(3, 1, 1, '"h, D, W, M, Q, Y"')

# --- Connections sheet ---------------------------------------------------
# Heat pipeline (pl_dh) now carries the low-grade heat node.
$connections.Range("C5").Value2 = "heat_low"
